$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.8546379068706642
$ws.Range("C2").Value = 0.1453813793327683
$ws.Range("D2").Value = 0.05331196988965559
$ws.Range("F2").Value = 1.248581532117257
$ws.Range("G2").Value = 1.124731381346066
$ws.Range("H2").Value = 1.099298929078799
$ws.Range("K2").Value = 0.4402374974239649
$ws.Range("L2").Value = 0.3041183138331007
$ws.Range("N2").Value = 2.054157730142066

# Row 3
$ws.Range("B3").Value = 0.8070385758574901
$ws.Range("C3").Value = 0.1447875784585975
$ws.Range("D3").Value = 0.05261386623917375
$ws.Range("F3").Value = 1.24167878236868
$ws.Range("G3").Value = 1.119222560094229
$ws.Range("H3").Value = 1.101473366892321
$ws.Range("K3").Value = 0.399407620758268
$ws.Range("L3").Value = 0.2931447837167553
$ws.Range("N3").Value = 2.072376104573351

# Row 4
$ws.Range("B4").Value = 0.7782410031614972
$ws.Range("C4").Value = 0.1444185682655075
$ws.Range("D4").Value = 0.05217782790180792
$ws.Range("F4").Value = 1.23812150568542
$ws.Range("G4").Value = 1.116483651846906
$ws.Range("H4").Value = 1.103276275730423
$ws.Range("K4").Value = 0.3744920737558459
$ws.Range("L4").Value = 0.2865851252705625
$ws.Range("N4").Value = 2.084209235229796

# Row 5
$ws.Range("B5").Value = 0.7666138606835773
$ws.Range("C5").Value = 0.1442670917370314
$ws.Range("D5").Value = 0.05199828836267528
$ws.Range("F5").Value = 1.236843069373904
$ws.Range("G5").Value = 1.115529180478745
$ws.Range("H5").Value = 1.104128640856047
$ws.Range("K5").Value = 0.3643777588673061
$ws.Range("L5").Value = 0.2839567887518228
$ws.Range("N5").Value = 2.089193930343693

# Row 6
$ws.Range("B6").Value = 0.7646897212974579
$ws.Range("C6").Value = 0.1442418728643773
$ws.Range("D6").Value = 0.05196836450886622
$ws.Range("F6").Value = 1.236641123017883
$ws.Range("G6").Value = 1.115380450074682
$ws.Range("H6").Value = 1.104277283158666
$ws.Range("K6").Value = 0.3627006446240841
$ws.Range("L6").Value = 0.2835230595036506
$ws.Range("N6").Value = 2.090031453150544

# Row 7
$ws.Range("B7").Value = 0.7780837573172619
$ws.Range("C7").Value = 0.1444165298506306
$ws.Range("D7").Value = 0.05217541404770643
$ws.Range("F7").Value = 1.238103571228983
$ws.Range("G7").Value = 1.11647012519434
$ws.Range("H7").Value = 1.103287294560445
$ws.Range("K7").Value = 0.3743555104922791
$ws.Range("L7").Value = 0.2865494973209763
$ws.Range("N7").Value = 2.084275802339935

# Row 8
$ws.Range("B8").Value = 0.8381368775703208
$ws.Range("C8").Value = 0.1451775575922838
$ws.Range("D8").Value = 0.05307280561254757
$ws.Range("F8").Value = 1.246060066249456
$ws.Range("G8").Value = 1.12269827079291
$ws.Range("H8").Value = 1.099951602706867
$ws.Range("K8").Value = 0.4261274705279448
$ws.Range("L8").Value = 0.3002976554409429
$ws.Range("N8").Value = 2.060305093828482

# Row 9
$ws.Range("B9").Value = 0.9592943975889909
$ws.Range("C9").Value = 0.1466346169324879
$ws.Range("D9").Value = 0.05477348885914779
$ws.Range("F9").Value = 1.267072178873434
$ws.Range("G9").Value = 1.140027412067354
$ws.Range("H9").Value = 1.097121502550536
$ws.Range("K9").Value = 0.5288716965494586
$ws.Range("L9").Value = 0.3286736726815889
$ws.Range("N9").Value = 2.018435489207832

# Row 10
$ws.Range("B10").Value = 1.0503781672586
$ws.Range("C10").Value = 0.1476832952272389
$ws.Range("D10").Value = 0.05598654196006692
$ws.Range("F10").Value = 1.285819096086016
$ws.Range("G10").Value = 1.155894367363146
$ws.Range("H10").Value = 1.097305302643832
$ws.Range("K10").Value = 0.6051045684659186
$ws.Range("L10").Value = 0.3503908061022969
$ws.Range("N10").Value = 1.990809377785887

# Row 11
$ws.Range("B11").Value = 1.092264656665009
$ws.Range("C11").Value = 0.1481555715254146
$ws.Range("D11").Value = 0.05653039967657492
$ws.Range("F11").Value = 1.295068882724252
$ws.Range("G11").Value = 1.163797218593686
$ws.Range("H11").Value = 1.097880583207839
$ws.Range("K11").Value = 0.6399482075570688
$ws.Range("L11").Value = 0.3604607173156893
$ws.Range("N11").Value = 1.978923254382728

# Row 12
$ws.Range("B12").Value = 1.108190814734428
$ws.Range("C12").Value = 0.1483337172038688
$ws.Range("D12").Value = 0.05673519005097205
$ws.Range("F12").Value = 1.298675470544779
$ws.Range("G12").Value = 1.166888555336755
$ws.Range("H12").Value = 1.098169129494025
$ws.Range("K12").Value = 0.6531662168677883
$ws.Range("L12").Value = 0.3643014122941679
$ws.Range("N12").Value = 1.974520345773165

# Row 13
$ws.Range("B13").Value = 1.10475796166935
$ws.Range("C13").Value = 0.1482953813989667
$ws.Range("D13").Value = 0.05669113640316681
$ws.Range("F13").Value = 1.297894104806872
$ws.Range("G13").Value = 1.166218387636704
$ws.Range("H13").Value = 1.098103841494435
$ws.Range("K13").Value = 0.6503184415673786
$ws.Range("L13").Value = 0.3634730292923933
$ws.Range("N13").Value = 1.975464224905046

# Row 14
$ws.Range("B14").Value = 1.093573617504433
$ws.Range("C14").Value = 0.1481702416591304
$ws.Range("D14").Value = 0.05654727115057057
$ws.Range("F14").Value = 1.295363516033561
$ws.Range("G14").Value = 1.164049565448494
$ws.Range("H14").Value = 1.097902905155109
$ws.Range("K14").Value = 0.6410351919282675
$ws.Range("L14").Value = 0.3607761435938102
$ws.Range("N14").Value = 1.97855905726334

# Row 15
$ws.Range("B15").Value = 1.08673129328298
$ws.Range("C15").Value = 0.1480934991724538
$ws.Range("D15").Value = 0.05645899855698389
$ws.Range("F15").Value = 1.293826990787409
$ws.Range("G15").Value = 1.162733958986692
$ws.Range("H15").Value = 1.09778903305363
$ws.Range("K15").Value = 0.6353519841804598
$ws.Range("L15").Value = 0.3591277970489415
$ws.Range("N15").Value = 1.980467512918786

# Row 16
$ws.Range("B16").Value = 1.047649852880738
$ws.Range("C16").Value = 0.1476523340679776
$ws.Range("D16").Value = 0.05595083831834557
$ws.Range("F16").Value = 1.285229133176955
$ws.Range("G16").Value = 1.155391697332036
$ws.Range("H16").Value = 1.097277600196819
$ws.Range("K16").Value = 0.6028307531657617
$ws.Range("L16").Value = 0.3497365522668048
$ws.Range("N16").Value = 1.991599883838738

# Row 17
$ws.Range("B17").Value = 1.023790182327673
$ws.Range("C17").Value = 0.1473804648761501
$ws.Range("D17").Value = 0.0556370501180723
$ws.Range("F17").Value = 1.280139555225247
$ws.Range("G17").Value = 1.151063025642998
$ws.Range("H17").Value = 1.097089779567185
$ws.Range("K17").Value = 0.5829221037472792
$ws.Range("L17").Value = 0.3440241784533526
$ws.Range("N17").Value = 1.998603801784292

# Row 18
$ws.Range("B18").Value = 1.010109298379064
$ws.Range("C18").Value = 0.147223644198931
$ws.Range("D18").Value = 0.05545581828746649
$ws.Range("F18").Value = 1.277280089554722
$ws.Range("G18").Value = 1.14863774089936
$ws.Range("H18").Value = 1.09702802831292
$ws.Range("K18").Value = 0.5714867126249032
$ws.Range("L18").Value = 0.3407565185437988
$ws.Range("N18").Value = 2.002696371693503

# Row 19
$ws.Range("B19").Value = 1.005484507747212
$ws.Range("C19").Value = 0.1471704706146895
$ws.Range("D19").Value = 0.05539432797869281
$ws.Range("F19").Value = 1.276323586393744
$ws.Range("G19").Value = 1.147827644544321
$ws.Range("H19").Value = 1.097015068792203
$ws.Range("K19").Value = 0.5676175614625265
$ws.Range("L19").Value = 0.3396532267849182
$ws.Range("N19").Value = 2.004093050168059

# Row 20
$ws.Range("B20").Value = 1.026325681933486
$ws.Range("C20").Value = 0.147409452319593
$ws.Range("D20").Value = 0.0556705310237291
$ws.Range("F20").Value = 1.280674318697464
$ws.Range("G20").Value = 1.151517147985771
$ws.Range("H20").Value = 1.09710498378341
$ws.Range("K20").Value = 0.5850398071116274
$ws.Range("L20").Value = 0.3446304128756594
$ws.Range("N20").Value = 1.997851587319872

# Row 21
$ws.Range("B21").Value = 1.096856978580149
$ws.Range("C21").Value = 0.1482070171690708
$ws.Range("D21").Value = 0.05658955934460863
$ws.Range("F21").Value = 1.296103990739766
$ws.Range("G21").Value = 1.164683920758549
$ws.Range("H21").Value = 1.097960006255107
$ws.Range("K21").Value = 0.64376127203343
$ws.Range("L21").Value = 0.3615675394829481
$ws.Range("N21").Value = 1.977647366037743

# Row 22
$ws.Range("B22").Value = 1.143329948457165
$ws.Range("C22").Value = 0.1487242162627425
$ws.Range("D22").Value = 0.05718344999873892
$ws.Range("F22").Value = 1.306793801719039
$ws.Range("G22").Value = 1.173864557083931
$ws.Range("H22").Value = 1.098930906926199
$ws.Range("K22").Value = 0.6822758905168769
$ws.Range("L22").Value = 0.3727968701223432
$ws.Range("N22").Value = 1.965014667683867

# Row 23
$ws.Range("B23").Value = 1.118492097055878
$ws.Range("C23").Value = 0.148448551350242
$ws.Range("D23").Value = 0.05686710044638232
$ws.Range("F23").Value = 1.301032996221366
$ws.Range("G23").Value = 1.168911960416139
$ws.Range("H23").Value = 1.098375012478158
$ws.Range("K23").Value = 0.6617074794153268
$ws.Range("L23").Value = 0.3667889243273521
$ws.Range("N23").Value = 1.971704590630971

# Row 24
$ws.Range("B24").Value = 1.025179269283854
$ws.Range("C24").Value = 0.1473963487146221
$ws.Range("D24").Value = 0.0556553968928668
$ws.Range("F24").Value = 1.280432344679284
$ws.Range("G24").Value = 1.151311642029242
$ws.Range("H24").Value = 1.097097965949189
$ws.Range("K24").Value = 0.5840823611460451
$ws.Range("L24").Value = 0.3443562830271105
$ws.Range("N24").Value = 1.998191458129234

# Row 25
$ws.Range("B25").Value = 0.9261547252329194
$ws.Range("C25").Value = 0.1462442551279075
$ws.Range("D25").Value = 0.05431977883119998
$ws.Range("F25").Value = 1.260807534048269
$ws.Range("G25").Value = 1.134790037636847
$ws.Range("H25").Value = 1.097489750949833
$ws.Range("K25").Value = 0.5009457691808166
$ws.Range("L25").Value = 0.320844994949141
$ws.Range("N25").Value = 2.029211997147272
